# Update cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.089.21'
$ws.Range("E2").Value = '  +0.64%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.682.02'
$ws.Range("E3").Value = '  +0.78%  '

# Row 4
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("E5").Value = '  +0.04%  '

# Row 6
$ws.Range("E6").Value = '  -3.46%  '

# Row 7
$ws.Range("E7").Value = '  +0.08%  '

# Row 8
$ws.Range("E8").Value = '  +1.65%  '

# Row 9
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0624'
$ws.Range("E9").Value = '  +0.50%  '

# Row 10
$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.30'
$ws.Range("E10").Value = '  +5.10%  '

# Row 11
$ws.Range("E11").Value = '  -0.58%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.919.38'
$ws.Range("E12").Value = '  +0.84%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.690.86'
$ws.Range("E13").Value = '  +1.40%  '

# Row 14
$ws.Range("E14").Value = '  +0.82%  '

# Row 15
$ws.Range("E15").Value = '  +1.99%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.44'
$ws.Range("E16").Value = '  +0.38%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.083.45'
$ws.Range("E17").Value = '  +0.53%  '

# Row 18
$ws.Range("E18").Value = '  +2.61%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '236.06'
$ws.Range("E19").Value = '  +0.66%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0737'
$ws.Range("E20").Value = '  +0.56%  '

# Row 21
$ws.Range("E21").Value = '  +0.01%  '

# Row 22
$ws.Range("E22").Value = '  +1.41%  '

# Row 23
$ws.Range("E23").Value = '  +1.29%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.14'
$ws.Range("E24").Value = '  -2.98%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.39'
$ws.Range("E25").Value = '  +0.89%  '

# Row 26
$ws.Range("E26").Value = '  +2.08%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.57'
$ws.Range("E27").Value = '  +4.04%  '

# Row 28
$ws.Range("E28").Value = '  -1.73%  '

# Row 29
$ws.Range("E29").Value = '  +0.21%  '

# Row 30
$ws.Range("E30").Value = '  +0.21%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.18'
$ws.Range("E31").Value = '  +0.46%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.38'
$ws.Range("E32").Value = '  +0.62%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.549.87'
$ws.Range("E33").Value = '  +6.73%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.17'
$ws.Range("E34").Value = '  +1.26%  '

# Row 35
$ws.Range("E35").Value = '  +4.16%  '

# Row 36
$ws.Range("E36").Value = '  -1.01%  '

# Row 37
$ws.Range("E37").Value = '  +0.74%  '

# Row 38
$ws.Range("E38").Value = '  +1.25%  '

# Row 39
$ws.Range("E39").Value = '  +2.48%  '

# Row 40
$ws.Range("E40").Value = '  +7.43%  '

# Row 41
$ws.Range("E41").Value = '  +0.06%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '67.94'
$ws.Range("E42").Value = '  +2.62%  '

# Row 43
$ws.Range("E43").Value = '  -3.50%  '

# Row 44
$ws.Range("E44").Value = '  -1.02%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.823.60'
$ws.Range("E45").Value = '  +0.80%  '

# Row 46
$ws.Range("E46").Value = '  -0.42%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.61'
$ws.Range("E47").Value = '  -0.01%  '

# Row 48
$ws.Range("E48").Value = '  +3.72%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.54'
$ws.Range("E49").Value = '  +0.34%  '

# Row 50
$ws.Range("E50").Value = '  +1.73%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.02'
$ws.Range("E51").Value = '  +6.52%  '

